$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.961.14"
$ws.Range("E2").Value = "  -4.60%  "
$ws.Range("D3").Value = "2.964.96"
$ws.Range("E3").Value = "  -4.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.91%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "2.961.48"
$ws.Range("E8").Value = "  -4.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.490"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("E10").Value = "  -7.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.87"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.440"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.61%  "
$ws.Range("E13").Value = "  -3.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.81%  "
$ws.Range("D15").Value = "3.446.66"
$ws.Range("E15").Value = "  -4.11%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "60.986.26"
$ws.Range("E16").Value = "  -4.71%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.109"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.30%  "
$ws.Range("D18").Value = "2.962.05"
$ws.Range("E18").Value = "  -4.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "462.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.659"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.82%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.56%  "
$ws.Range("E32").Value = "  -4.37%  "
$ws.Range("E33").Value = "  -4.39%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "54.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.02%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "442.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -11.05%  "
$ws.Range("D38").Value = "3.121.42"
$ws.Range("E38").Value = "  -3.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0779"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.57%  "
$ws.Range("E40").Value = "  -3.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0373"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("E44").Value = "  -13.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("E46").Value = "  -6.42%  "
$ws.Range("E47").Value = "  -2.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "115.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.49%  "
$ws.Range("E51").Value = "  -11.04%  "
